# Appends the (currently empty) data row #2 under the header row of the
# "Don sale chinh" report sheet.
#
# The report is generated/refreshed by a per-employee export routine that
# always writes one row per record across a fixed 20-column (A:T) schema:
# text columns default to an empty string, money/quantity columns default
# to 0. For this employee there is no underlying data yet, so row 2 comes
# out "empty" - but it is still materialised (sheet dimension grows to
# A1:T2).
#
# Money / quantity columns (original price, up-sale, unit price, first
# payment, pay-later, paid, debt) are zero-filled; every text column
# (prefix, service code, dates, branch, customer, source, service name,
# sales reps, doctors, assistants) is cleared to an empty value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Đơn sale chính")

$row = 2

$textCols    = @("A", "B", "C", "D", "E", "F", "G", "H", "J", "Q", "R", "S", "T")
$numZeroCols = @("I", "K", "L", "M", "N", "O", "P")

foreach ($col in $textCols) {
    $ws.Range("$col$row").Value = ""
}

foreach ($col in $numZeroCols) {
    $ws.Range("$col$row").Value = 0
}
